$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was added to the underlying ifoCAST full-series
# evaluation. This shifts every rolling-window statistic (columns B:G,
# rows 2-11) up by one row, and row 11 (the last row) now receives the
# newly computed statistics for the latest window.

$ws.Range("B2").Value = -0.06608967275348025
$ws.Range("C2").Value = 0.7088863950418706
$ws.Range("D2").Value = 1.217845040002488
$ws.Range("E2").Value = 1.103560166009306
$ws.Range("F2").Value = 1.133515860779062
$ws.Range("G2").Value = 18

$ws.Range("B3").Value = 0.001870457576035796
$ws.Range("C3").Value = 0.6187918867460556
$ws.Range("D3").Value = 0.6230069128923086
$ws.Range("E3").Value = 0.7893078695238688
$ws.Range("F3").Value = 0.8135976448302396
$ws.Range("G3").Value = 17

$ws.Range("B4").Value = 0.1204196123905843
$ws.Range("C4").Value = 0.5277616828532142
$ws.Range("D4").Value = 0.6040816134619066
$ws.Range("E4").Value = 0.7772268738675385
$ws.Range("F4").Value = 0.7930233996807294
$ws.Range("G4").Value = 16

$ws.Range("B5").Value = 0.09828742272023362
$ws.Range("C5").Value = 0.5953749172355302
$ws.Range("D5").Value = 0.77910483569946
$ws.Range("E5").Value = 0.8826691541565617
$ws.Range("F5").Value = 0.9079673638532668
$ws.Range("G5").Value = 15

$ws.Range("B6").Value = 0.1268913055883079
$ws.Range("C6").Value = 0.6888184636477384
$ws.Range("D6").Value = 1.036376752836633
$ws.Range("E6").Value = 1.018025909707917
$ws.Range("F6").Value = 1.048216613390359
$ws.Range("G6").Value = 14

$ws.Range("B7").Value = 0.1309701764864839
$ws.Range("C7").Value = 0.6044880343376112
$ws.Range("D7").Value = 0.7446009131626561
$ws.Range("E7").Value = 0.8629026093150119
$ws.Range("F7").Value = 0.8877321498458313
$ws.Range("G7").Value = 13

$ws.Range("B8").Value = 0.1991520669947553
$ws.Range("C8").Value = 0.6396550887298748
$ws.Range("D8").Value = 0.7959921877658299
$ws.Range("E8").Value = 0.8921839427863684
$ws.Range("F8").Value = 0.9083435325175242
$ws.Range("G8").Value = 12

$ws.Range("B9").Value = 0.1894669389442902
$ws.Range("C9").Value = 0.7879769448546992
$ws.Range("D9").Value = 1.119193970586897
$ws.Range("E9").Value = 1.057919642783372
$ws.Range("F9").Value = 1.09161617549273
$ws.Range("G9").Value = 11

$ws.Range("B10").Value = 0.1730349822307496
$ws.Range("C10").Value = 0.6656591560163745
$ws.Range("D10").Value = 0.9145633670025062
$ws.Range("E10").Value = 0.9563280645272867
$ws.Range("F10").Value = 0.9914200040160747
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.2315707031541324
$ws.Range("C11").Value = 0.8270225223868328
$ws.Range("D11").Value = 0.9737462206433012
$ws.Range("E11").Value = 0.9867858028180692
$ws.Range("F11").Value = 1.01741652426354
$ws.Range("G11").Value = 9
